$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Anxa2"
$ws.Cells.Item(2,3).Value = "Tlr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 120.6115993333333
$ws.Cells.Item(2,8).Value = 361.834798
$ws.Cells.Item(2,9).Value = 0.4273073648704228
$ws.Cells.Item(2,10).Value = 0.4273073648704228
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 57.24915866666667
$ws.Cells.Item(2,14).Value = 171.747476
$ws.Cells.Item(2,15).Value = 0.9704198736548433
$ws.Cells.Item(2,16).Value = 0.9704198736548435
$ws.Cells.Item(2,17).Value = 6904.912587274427
$ws.Cells.Item(2,18).Value = 62144.21328546984
$ws.Cells.Item(2,19).Value = 0.4146675590293397
$ws.Cells.Item(2,20).Value = 0.4146675590293398

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Anxa2"
$ws.Cells.Item(3,3).Value = "Tlr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 120.6115993333333
$ws.Cells.Item(3,8).Value = 361.834798
$ws.Cells.Item(3,9).Value = 0.4273073648704228
$ws.Cells.Item(3,10).Value = 0.4273073648704228
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.569166
$ws.Cells.Item(3,14).Value = 4.707498
$ws.Cells.Item(3,15).Value = 0.02659864191768634
$ws.Cells.Item(3,16).Value = 0.02659864191768634
$ws.Cells.Item(3,17).Value = 189.2596208794893
$ws.Cells.Item(3,18).Value = 1703.336587915404
$ws.Cells.Item(3,19).Value = 0.01136579558697852
$ws.Cells.Item(3,20).Value = 0.01136579558697852

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Anxa2"
$ws.Cells.Item(4,3).Value = "Tlr2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 120.6115993333333
$ws.Cells.Item(4,8).Value = 361.834798
$ws.Cells.Item(4,9).Value = 0.4273073648704228
$ws.Cells.Item(4,10).Value = 0.4273073648704228
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.1758903333333333
$ws.Cells.Item(4,14).Value = 0.527671
$ws.Cells.Item(4,15).Value = 0.002981484427470275
$ws.Cells.Item(4,16).Value = 0.002981484427470276
$ws.Cells.Item(4,17).Value = 21.21441441060644
$ws.Cells.Item(4,18).Value = 190.929729695458
$ws.Cells.Item(4,19).Value = 0.001274010254104524
$ws.Cells.Item(4,20).Value = 0.001274010254104525

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Anxa2"
$ws.Cells.Item(5,3).Value = "Tlr2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 134.43072
$ws.Cells.Item(5,8).Value = 403.29216
$ws.Cells.Item(5,9).Value = 0.4762662715555095
$ws.Cells.Item(5,10).Value = 0.4762662715555095
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 57.24915866666667
$ws.Cells.Item(5,14).Value = 171.747476
$ws.Cells.Item(5,15).Value = 0.9704198736548433
$ws.Cells.Item(5,16).Value = 0.9704198736548435
$ws.Cells.Item(5,17).Value = 7696.04561895424
$ws.Cells.Item(5,18).Value = 69264.41057058817
$ws.Cells.Item(5,19).Value = 0.4621782550689609
$ws.Cells.Item(5,20).Value = 0.462178255068961

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Anxa2"
$ws.Cells.Item(6,3).Value = "Tlr2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 134.43072
$ws.Cells.Item(6,8).Value = 403.29216
$ws.Cells.Item(6,9).Value = 0.4762662715555095
$ws.Cells.Item(6,10).Value = 0.4762662715555095
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.569166
$ws.Cells.Item(6,14).Value = 4.707498
$ws.Cells.Item(6,15).Value = 0.02659864191768634
$ws.Cells.Item(6,16).Value = 0.02659864191768634
$ws.Cells.Item(6,17).Value = 210.94411517952
$ws.Cells.Item(6,18).Value = 1898.49703661568
$ws.Cells.Item(6,19).Value = 0.01266803601457656
$ws.Cells.Item(6,20).Value = 0.01266803601457656

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Anxa2"
$ws.Cells.Item(7,3).Value = "Tlr2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 134.43072
$ws.Cells.Item(7,8).Value = 403.29216
$ws.Cells.Item(7,9).Value = 0.4762662715555095
$ws.Cells.Item(7,10).Value = 0.4762662715555095
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.1758903333333333
$ws.Cells.Item(7,14).Value = 0.527671
$ws.Cells.Item(7,15).Value = 0.002981484427470275
$ws.Cells.Item(7,16).Value = 0.002981484427470276
$ws.Cells.Item(7,17).Value = 23.64506415104
$ws.Cells.Item(7,18).Value = 212.80557735936
$ws.Cells.Item(7,19).Value = 0.001419980471972081
$ws.Cells.Item(7,20).Value = 0.001419980471972081

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Anxa2"
$ws.Cells.Item(8,3).Value = "Tlr2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 27.21726533333333
$ws.Cells.Item(8,8).Value = 81.651796
$ws.Cells.Item(8,9).Value = 0.09642636357406766
$ws.Cells.Item(8,10).Value = 0.09642636357406766
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 57.24915866666667
$ws.Cells.Item(8,14).Value = 171.747476
$ws.Cells.Item(8,15).Value = 0.9704198736548433
$ws.Cells.Item(8,16).Value = 0.9704198736548435
$ws.Cells.Item(8,17).Value = 1558.165541540766
$ws.Cells.Item(8,18).Value = 14023.4898738669
$ws.Cells.Item(8,19).Value = 0.09357405955654272
$ws.Cells.Item(8,20).Value = 0.09357405955654274

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Anxa2"
$ws.Cells.Item(9,3).Value = "Tlr2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 27.21726533333333
$ws.Cells.Item(9,8).Value = 81.651796
$ws.Cells.Item(9,9).Value = 0.09642636357406766
$ws.Cells.Item(9,10).Value = 0.09642636357406766
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.569166
$ws.Cells.Item(9,14).Value = 4.707498
$ws.Cells.Item(9,15).Value = 0.02659864191768634
$ws.Cells.Item(9,16).Value = 0.02659864191768634
$ws.Cells.Item(9,17).Value = 42.70840737404534
$ws.Cells.Item(9,18).Value = 384.375666366408
$ws.Cells.Item(9,19).Value = 0.002564810316131259
$ws.Cells.Item(9,20).Value = 0.002564810316131259

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Anxa2"
$ws.Cells.Item(10,3).Value = "Tlr2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 27.21726533333333
$ws.Cells.Item(10,8).Value = 81.651796
$ws.Cells.Item(10,9).Value = 0.09642636357406766
$ws.Cells.Item(10,10).Value = 0.09642636357406766
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.1758903333333333
$ws.Cells.Item(10,14).Value = 0.527671
$ws.Cells.Item(10,15).Value = 0.002981484427470275
$ws.Cells.Item(10,16).Value = 0.002981484427470276
$ws.Cells.Item(10,17).Value = 4.787253871901778
$ws.Cells.Item(10,18).Value = 43.085284847116
$ws.Cells.Item(10,19).Value = 0.0002874937013936697
$ws.Cells.Item(10,20).Value = 0.0002874937013936698
